# Updates the cryptos list (prices / 1h volume %) per the Sun Oct  8 06:09:03 UTC 2023
# GitHub Actions refresh. Rows 33/34 (InternetComputer <-> Maker) and 49-51
# (BabyDogeCoin dropped, Algorand/Cronos shift up, EnergySwap appended) are
# full row replacements; the rest are in-place price/volume refreshes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.955.01'
$ws.Range('E2').Value = '  +0.18%  '
$ws.Range('D3').Value = '1.630.79'
$ws.Range('E3').Value = '  -0.58%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '211.86'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.88%  '
$ws.Range('E6').Value = '  -0.08%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '23.43'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -1.35%  '
$ws.Range('E9').Value = '  -2.09%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0614'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -0.29%  '
$ws.Range('E11').Value = '  +0.62%  '
$ws.Range('D12').Value = '1.861.58'
$ws.Range('E12').Value = '  -0.64%  '
$ws.Range('D13').Value = '1.625.21'
$ws.Range('E13').Value = '  -0.94%  '
$ws.Range('E14').Value = '  -1.17%  '
$ws.Range('E15').Value = '  -2.51%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '65.60'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -0.85%  '
$ws.Range('D17').Value = '27.943.52'
$ws.Range('E17').Value = '  +0.14%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '230.97'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -0.42%  '
$ws.Range('E19').Value = '  -0.09%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.64'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.37%  '
$ws.Range('E21').Value = '  -0.06%  '
$ws.Range('E22').Value = '  -9.60%  '
$ws.Range('E23').Value = '  -0.88%  '
$ws.Range('E24').Value = '  -1.35%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '155.32'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +2.23%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '6.95'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('E27').Value = '  -0.39%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '15.59'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -1.03%  '
$ws.Range('E29').Value = '  +0.03%  '
$ws.Range('E30').Value = '  -0.57%  '
$ws.Range('E31').Value = '  -0.79%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.40'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +1.63%  '
$ws.Range('B33').Value = 'Maker'
$ws.Range('C33').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D33').Value = '1.401.31'
$ws.Range('E33').Value = '  -1.57%  '
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.07'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -1.60%  '
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('E36').Value = '  +12.45%  '
$ws.Range('E38').Value = '  +2.02%  '
$ws.Range('E39').Value = '  -0.48%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.866'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -2.99%  '
$ws.Range('E41').Value = '  -0.67%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '66.60'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.99%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.82'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('E45').Value = '  -0.14%  '
$ws.Range('E46').Value = '  -0.41%  '
$ws.Range('D47').Value = '1.771.76'
$ws.Range('E47').Value = '  -0.60%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '88.20'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -0.47%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.1000'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.70%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0504'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -0.41%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '7.55'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -1.85%  '
